# Electricity Technology Shareweights - calibrate ETS (onshore wind row)
$wb = $excel.ActiveWorkbook
$wsETS = $wb.Worksheets.Item("ETS")

# Row 6 = "onshore wind" shareweights.
# 2020-2030 (cols B:L) -> 3
# 2031-2040 (cols M:V) -> 7
# 2041-2050 (cols W:AF) -> 10
$wsETS.Range("B6:L6").Value = 3
$wsETS.Range("M6:V6").Value = 7
$wsETS.Range("W6:AF6").Value = 10

# The header cell A1 is no longer italicized.
$wsETS.Range("A1").Font.Italic = $false

# Make the ETS sheet the active/selected tab, with A1:AF17 selected.
$wsETS.Activate() | Out-Null
$wsETS.Range("A1:AF17").Select() | Out-Null
